$wb = $excel.ActiveWorkbook

# --- Update the conversion note text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.76 = 6529.16 pesos`n✅ 6529.16 pesos = 1.74 = 932.49 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update rate values on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 568.9
$wsTasas.Range("O10").Value = 3714.44
$wsTasas.Range("N12").Value = 3746
$wsTasas.Range("O12").Value = 535
